$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing cell O32 (row 32, column 15) from 0 to 242
$ws.Cells.Item(32, 15).Value = 242

# Add new row 33 with the TCHD update values
$ws.Cells.Item(33, 1).Value = 31
# Copy the formatting (bold, bordered, centered style) from A32 onto A33
$ws.Range("A32").Copy()
$ws.Range("A33").PasteSpecial(-4122)

$ws.Cells.Item(33, 2).Value = 0
$ws.Cells.Item(33, 3).Value = 0
$ws.Cells.Item(33, 4).Value = 0
$ws.Cells.Item(33, 5).Value = 0
$ws.Cells.Item(33, 6).Value = 0
$ws.Cells.Item(33, 7).Value = 0
$ws.Cells.Item(33, 8).Value = 0
$ws.Cells.Item(33, 9).Value = 0
$ws.Cells.Item(33, 10).Value = 0
$ws.Cells.Item(33, 11).Value = 0
$ws.Cells.Item(33, 12).Value = 0
$ws.Cells.Item(33, 13).Value = 0
$ws.Cells.Item(33, 14).Value = 0
$ws.Cells.Item(33, 15).Value = 242
$ws.Cells.Item(33, 16).Value = 80
$ws.Cells.Item(33, 17).Value = 1158
$ws.Cells.Item(33, 18).Value = 1480
$ws.Cells.Item(33, 19).Value = 0
$ws.Cells.Item(33, 20).Value = 0
$ws.Cells.Item(33, 21).Value = 0
$ws.Cells.Item(33, 22).Value = 1
$ws.Cells.Item(33, 23).Value = 1
$ws.Cells.Item(33, 24).Value = 0
$ws.Cells.Item(33, 25).Value = 33
